$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2677.8096
$ws.Range("I51").Value = 2200
$ws.Range("J51").Value = 2827.125
$ws.Range("K51").Value = 2200
$ws.Range("L51").Value = 2827.125
$ws.Range("M51").Value = -1716
$ws.Range("N51").Value = -3795.125
# Row 53
$ws.Range("H53").Value = 1664.7858
$ws.Range("I53").Value = 2289.6
$ws.Range("J53").Value = 102.75
$ws.Range("K53").Value = 2289.6
$ws.Range("L53").Value = 102.75
$ws.Range("M53").Value = -1652.6
$ws.Range("N53").Value = -1376.75
# Row 100
$ws.Range("H100").Value = 3365.875
$ws.Range("I100").Value = 2742.3333
$ws.Range("J100").Value = 3740
$ws.Range("K100").Value = 2742.3333
$ws.Range("L100").Value = 3740
$ws.Range("M100").Value = -2201.3333
$ws.Range("N100").Value = -4822
# Row 129
$ws.Range("H129").Value = 2457.5686
$ws.Range("I129").Value = 485.5
$ws.Range("J129").Value = 3064.359
$ws.Range("K129").Value = 1456.5
$ws.Range("L129").Value = 9193.076999999999
$ws.Range("M129").Value = 3543.5
$ws.Range("N129").Value = -19193.077
# Row 138
$ws.Range("H138").Value = 6945795
$ws.Range("I138").Value = 1216.4615
$ws.Range("J138").Value = 25001698
$ws.Range("K138").Value = 3649.3845
$ws.Range("L138").Value = 75005094
$ws.Range("M138").Value = 1490.6155
$ws.Range("N138").Value = -75015374

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# Row 61
$ws.Range("H61").Value = 17243986
$ws.Range("I61").Value = 17243986
$ws.Range("K61").Value = 17243986
$ws.Range("M61").Value = -17243774
# Row 122
$ws.Range("H122").Value = 5912.615
$ws.Range("I122").Value = 5912.615
$ws.Range("K122").Value = 17737.845
$ws.Range("M122").Value = -15287.845
# Row 132
$ws.Range("H132").Value = 8931142
$ws.Range("I132").Value = 13891241
$ws.Range("J132").Value = 2962.1
$ws.Range("K132").Value = 41673723
$ws.Range("L132").Value = 8886.299999999999
$ws.Range("M132").Value = -41671193
$ws.Range("N132").Value = -13946.3
# Row 136
$ws.Range("H136").Value = 17243986
$ws.Range("I136").Value = 17243986
$ws.Range("K136").Value = 51731958
$ws.Range("M136").Value = -51729408

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1207.8235
$ws.Range("I20").Value = 1450.6666
$ws.Range("J20").Value = 934.625
$ws.Range("K20").Value = 1450.6666
$ws.Range("L20").Value = 934.625
$ws.Range("M20").Value = -1203.6666
$ws.Range("N20").Value = -1428.625
# Row 134
$ws.Range("H134").Value = 8283.5
$ws.Range("I134").Value = 7776
$ws.Range("J134").Value = 9400
$ws.Range("K134").Value = 23328
$ws.Range("L134").Value = 28200
$ws.Range("M134").Value = -20793
$ws.Range("N134").Value = -33270

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 83338216
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 83338216
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 83338216
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -83338806
# Row 34
$ws.Range("H34").Value = 83338216
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 83338216
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 83338216
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -83338620
# Row 99
$ws.Range("H99").Value = 1131
$ws.Range("I99").Value = 1061.2858
$ws.Range("J99").Value = 1375
$ws.Range("K99").Value = 1061.2858
$ws.Range("L99").Value = 1375
$ws.Range("M99").Value = 436.7141999999999
$ws.Range("N99").Value = -4371
# Row 126
$ws.Range("H126").Value = 1131
$ws.Range("I126").Value = 1061.2858
$ws.Range("J126").Value = 1375
$ws.Range("K126").Value = 3183.8574
$ws.Range("L126").Value = 4125
$ws.Range("M126").Value = -713.8574000000003
$ws.Range("N126").Value = -9065
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 939.3
$ws.Range("J86").Value = 1032.1666
$ws.Range("L86").Value = 3096.4998
$ws.Range("N86").Value = -5468.4998
# Row 89
$ws.Range("H89").Value = 939.3
$ws.Range("J89").Value = 1032.1666
$ws.Range("L89").Value = 9289.499400000001
$ws.Range("N89").Value = -21145.4994
# Row 98
$ws.Range("H98").Value = 303.93332
$ws.Range("J98").Value = 187.5
$ws.Range("L98").Value = 562.5
$ws.Range("N98").Value = -3558.5
# Row 125
$ws.Range("H125").Value = 5620
$ws.Range("J125").Value = 5620
$ws.Range("L125").Value = 16860
$ws.Range("N125").Value = -26700
# Row 129
$ws.Range("H129").Value = 4269.9
$ws.Range("I129").Value = 1350
$ws.Range("J129").Value = 8649.75
$ws.Range("K129").Value = 4050
$ws.Range("L129").Value = 25949.25
$ws.Range("M129").Value = 950
$ws.Range("N129").Value = -35949.25
# Row 130
$ws.Range("H130").Value = 5662.143
$ws.Range("I130").Value = 1054
$ws.Range("J130").Value = 8222.223
$ws.Range("K130").Value = 3162
$ws.Range("L130").Value = 24666.669
$ws.Range("M130").Value = 1858
$ws.Range("N130").Value = -34706.669
# Row 131
$ws.Range("H131").Value = 3118303
$ws.Range("I131").Value = 7335784.5
$ws.Range("J131").Value = 1033.9131
$ws.Range("K131").Value = 22007353.5
$ws.Range("L131").Value = 3101.7393
$ws.Range("M131").Value = -22002313.5
$ws.Range("N131").Value = -13181.7393

$ws = $wb.Worksheets.Item("GSM")
# Row 130
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 5376.087
$ws.Range("I132").Value = 4410
$ws.Range("J132").Value = 6261.6665
$ws.Range("K132").Value = 13230
$ws.Range("L132").Value = 18784.9995
$ws.Range("M132").Value = -10700
$ws.Range("N132").Value = -23844.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6673.6
$ws.Range("I7").Value = 6960.4
$ws.Range("J7").Value = 6100
$ws.Range("K7").Value = 6960.4
$ws.Range("L7").Value = 6100
$ws.Range("M7").Value = -6848.4
$ws.Range("N7").Value = -6324
# Row 22
$ws.Range("H22").Value = 1114.1034
$ws.Range("I22").Value = 803.5625
$ws.Range("K22").Value = 803.5625
$ws.Range("M22").Value = -508.5625
# Row 27
$ws.Range("H27").Value = 1114.1034
$ws.Range("I27").Value = 803.5625
$ws.Range("K27").Value = 803.5625
$ws.Range("M27").Value = -696.5625
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 126
$ws.Range("H126").Value = 6673.6
$ws.Range("I126").Value = 6960.4
$ws.Range("J126").Value = 6100
$ws.Range("K126").Value = 20881.2
$ws.Range("L126").Value = 18300
$ws.Range("M126").Value = -18411.2
$ws.Range("N126").Value = -23240

$ws = $wb.Worksheets.Item("WVR")
# Row 124
$ws.Range("H124").Value = 56500
$ws.Range("J124").Value = 56500
$ws.Range("L124").Value = 56500
$ws.Range("N124").Value = -66320
# Row 126
$ws.Range("H126").Value = 2738.3928
$ws.Range("I126").Value = 1615
$ws.Range("J126").Value = 12100
$ws.Range("K126").Value = 4845
$ws.Range("L126").Value = 36300
$ws.Range("M126").Value = -2375
$ws.Range("N126").Value = -41240
# Row 127
$ws.Range("H127").Value = 26500
$ws.Range("I127").Value = 26500
$ws.Range("K127").Value = 26500
$ws.Range("M127").Value = -21540
